$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date) and Volumen values between row 2 and row 5
$ws.Range("D2").Value = 44257
$ws.Range("M2").Value = 100

$ws.Range("D5").Value = 44250
$ws.Range("M5").Value = 200
